# 自动更新Excel文件 - 2026-01-21 23:18:22
#
# For every data row (2..99) on the active sheet:
#   D = total days, E = remaining days, F = start date (yyyyMMdd)
# The "current date" advances from 2026-01-21 to 2026-01-22.
# Remaining days = (start date + total days) - current date.
# If that would be <= 0 (i.e. the period has lapsed as of the new date),
# the record is renewed: the start date becomes the new current date and
# the remaining days reset back to the total days.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$todayStr = "20260122"
$today = [datetime]::ParseExact($todayStr, "yyyyMMdd", $null)

$lastRow = 99

for ($row = 2; $row -le $lastRow; $row++) {
    $dCell = $ws.Range("D$row")
    $fCell = $ws.Range("F$row")
    $eCell = $ws.Range("E$row")

    $d = $dCell.Value2
    $fRaw = $fCell.Value2

    if ($d -eq $null -or $fRaw -eq $null) {
        continue
    }

    $fStr = [string]$fRaw

    try {
        $start = [datetime]::ParseExact($fStr, "yyyyMMdd", $null)
    } catch {
        # Unparseable / malformed start date (e.g. data typo) - leave row untouched.
        continue
    }

    $end = $start.AddDays($d)
    $remaining = $end.ToOADate() - $today.ToOADate()

    if ($remaining -le 0) {
        # Period lapsed - renew from today.
        $remaining = $d
        $fStr = $todayStr
    }

    $eCell.Value = $remaining
    $fCell.Value = [int]$fStr
}
